# Pulumi environment set up
# Replace the placeholder "Tenant" row with real Pulumi/Azure configuration
# values, and add a new "Pulumi Storage Account" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: Subscription label / value (overwrites old "Tenant" / "xxx")
$ws.Range("A1").Value = "Subscription"
$ws.Range("B1").Value = "Azure subscription 1"

# Row 2: Pulumi Resource Group label / value (overwrites old "Subscription" / "yyy")
$ws.Range("A2").Value = "Pulumi Resource Group"
$ws.Range("B2").Value = "rg-pulumi-nvtst"

# Row 3: new Pulumi Storage Account label / value (replaces old "zzz")
$ws.Range("A3").Value = "Pulumi Storage Account"
$ws.Range("B3").Value = "stpuluminvtst"

# Widen column B to fit the new values and select the whole of row 2
$ws.Columns.Item(2).ColumnWidth = 26.5
[void]$ws.Range("A2:XFD2").Select()
